# edit.ps1
# Applies two changes described by the diff:
#  1. In the "Ex 3" paragraph, merge the two runs that were split by a
#     "_GoBack" bookmark ("...pode ser " + bookmark + "planejada...pacotes.")
#     into a single run, removing the bookmark from this location.
#  2. Fill in the following (previously empty) paragraph with the full
#     "Ex 4" text, including its proofErr spell/gram-check markers and a
#     "_GoBack" bookmark in its new location.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge "pode ser " / bookmark / "planejada ... pacotes." runs
# ---------------------------------------------------------------------------

# Find the "Ex 3" paragraph by locating the SCRUM text, then resolve the
# enclosing paragraph so the edit is not dependent on a fixed paragraph index.
$locate = $d.Content
$locate.Find.Execute("SCRUM") | Out-Null
$ex3 = $locate.Paragraphs(1)

# The replaced span starts right after ". O " and ends right after "pacotes."
# (exclusive of the paragraph mark), i.e. it covers the two runs that used to
# be split by the bookmark, but leaves the ". O " run (and its proofErr
# markers) intact.
$beginRange = $d.Range($ex3.Range.Start, $ex3.Range.End)
$beginRange.Find.Execute(". O ") | Out-Null
$beginPos = $beginRange.End

$endRange = $d.Range($ex3.Range.Start, $ex3.Range.End)
$endRange.Find.Execute("pacotes.") | Out-Null
$endPos = $endRange.End

$mergedRunXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>método SCRUM pode ser aplicado por ser um método Ágil de produção, e, portanto, pode ser planejada a entrega das tarefas, tendo o status diário delas em uma reunião rotineira, podendo ser entregue versões do projeto para o cliente em pacotes.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target1 = $d.Range($beginPos, $endPos)
$target1.InsertXML($mergedRunXml)

# ---------------------------------------------------------------------------
# Change 2: fill the empty paragraph after "Ex 3" with the "Ex 4" content
# ---------------------------------------------------------------------------

# The target is the (now still empty) paragraph that immediately follows the
# "Ex 3" paragraph edited above.
$ex4Para = $ex3.Next()

$ex4Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>Ex</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>4</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve">. A equipe pode ser dividida por um </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>PO(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>Product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>Owner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve">) que deve conhecer tudo sobre o produto que será entregue, o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>Scrum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve"> Master que tem como o principal objetivo capacitar o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>squad</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve"> e auxiliá-los caso tenha algo pendente de ou</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve">tro setor, e o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>squad</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve">, que pode ser composto por </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>full</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>stack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>developers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve"> que irão ser responsáveis pelo desenvolvimento do software na área de front-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>end</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t>back-end</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr><w:t xml:space="preserve"> e na administração de banco de dados.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# Insert at a collapsed range located at the paragraph-mark position (the
# paragraph's End): this replaces the (empty) paragraph's content in place
# while preserving its original <w:p>/<w:pPr> attributes, instead of pushing
# a new empty paragraph after it.
$insertPos = $ex4Para.Range.End
$target2 = $d.Range($insertPos, $insertPos)
$target2.InsertXML($ex4Xml)
